# Update "想去人数" (wanted-to-go count) figures on the 展览 and 全部类型 sheets.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F15").Value = 4362
$ws1.Range("F16").Value = 1267
$ws1.Range("F18").Value = 2757
$ws1.Range("F20").Value = 12
$ws1.Range("F22").Value = 3769
$ws1.Range("F31").Value = 989
$ws1.Range("F44").Value = 612
$ws1.Range("F48").Value = 250

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F17").Value = 4362
$ws4.Range("F18").Value = 1267
$ws4.Range("F21").Value = 2757
$ws4.Range("F23").Value = 3769
$ws4.Range("F36").Value = 989
$ws4.Range("F44").Value = 612
$ws4.Range("F48").Value = 250
